$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Artfynd rows 5-17 were reshuffled (and row 7 gained a "K" Ålder-Stadium note
# that used to sit on row 14) per the upstream re-export. Columns: A=Id,
# B=Taxonsorteringsordning, D=Rödlistade, E=TaxonId, F=Artnamn,
# G=Vetenskapligt namn, H=Auktor, K=Ålder-Stadium, Q=Ost, R=Nord.
$rows = @(
    @{ Row = 5; A = 111943990; B = 101703; D = 'LC'; E = 222412; F = 'Tibast'; G = 'Daphne mezereum'; H = 'L.'; K = $null; Q = 682930.0967543643; R = 6694720.015570021 }
    @{ Row = 6; A = 111943980; B = 89183; D = 'LC'; E = 3215; F = 'Rödgul trumpetsvamp'; G = 'Craterellus lutescens'; H = '(Fr.) Fr.'; K = $null; Q = 682877.1417635784; R = 6694410.432217407 }
    @{ Row = 7; A = 111943997; B = 96326; D = 'LC'; E = 219798; F = 'Skogsknipprot'; G = 'Epipactis helleborine'; H = '(L.) Crantz'; K = 'i frukt'; Q = 682780.8405377725; R = 6694488.393080305 }
    @{ Row = 8; A = 111943995; B = 88899; D = 'NT'; E = 3286; F = 'Flattoppad klubbsvamp'; G = 'Clavariadelphus truncatus'; H = '(Quél.) Donk'; K = $null; Q = 682779.1674098044; R = 6694551.279700429 }
    @{ Row = 9; A = 111943988; B = 107033; D = 'NT'; E = 220320; F = 'Ängsskära'; G = 'Serratula tinctoria'; H = 'L.'; K = $null; Q = 682930.0967543643; R = 6694720.015570021 }
    @{ Row = 10; A = 111943984; B = 99413; D = 'LC'; E = 221235; F = 'Vårärt'; G = 'Lathyrus vernus'; H = '(L.) Bernh.'; K = $null; Q = 682929.3627028114; R = 6694685.271877083 }
    @{ Row = 11; A = 111943999; B = 99413; D = 'LC'; E = 221235; F = 'Vårärt'; G = 'Lathyrus vernus'; H = '(L.) Bernh.'; K = $null; Q = 682757.1772001419; R = 6694405.884787144 }
    @{ Row = 12; A = 111943996; B = 90332; D = 'LC'; E = 4769; F = 'Svavelriska'; G = 'Lactarius scrobiculatus'; H = '(Scop.:Fr.) Fr.'; K = $null; Q = 682785.3360249697; R = 6694547.127516991 }
    @{ Row = 13; A = 111943981; B = 96253; D = 'LC'; E = 504; F = 'Guckusko'; G = 'Cypripedium calceolus'; H = 'L.'; K = $null; Q = 682877.1417635784; R = 6694410.432217407 }
    @{ Row = 14; A = 111943979; B = 96253; D = 'LC'; E = 504; F = 'Guckusko'; G = 'Cypripedium calceolus'; H = 'L.'; K = $null; Q = 682878.8271195606; R = 6694406.550233844 }
    @{ Row = 15; A = 111943998; B = 98535; D = 'LC'; E = 222498; F = 'Blåsippa'; G = 'Hepatica nobilis'; H = 'Schreb.'; K = $null; Q = 682757.1772001419; R = 6694405.884787144 }
    @{ Row = 16; A = 111943983; B = 90678; D = 'LC'; E = 4366; F = 'Skarp dropptaggsvamp'; G = 'Hydnellum peckii'; H = 'Banker'; K = $null; Q = 682871.1304590552; R = 6694480.539619928 }
    @{ Row = 17; A = 111943992; B = 89183; D = 'LC'; E = 3215; F = 'Rödgul trumpetsvamp'; G = 'Craterellus lutescens'; H = '(Fr.) Fr.'; K = $null; Q = 682866.8554180798; R = 6694644.443727687 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    if ($null -eq $r.K) {
        $ws.Cells.Item($r.Row, 11).Value = ""
    } else {
        $ws.Cells.Item($r.Row, 11).Value = $r.K
    }
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
    $ws.Cells.Item($r.Row, 18).Value = $r.R
}
